$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap E10 / F10 (Thursday / Friday lab on row 10)
$ws.Range("E10").Value = "HT LAB"
$ws.Range("F10").Value = "TOM LAB"

# Swap E13 / F13 (Thursday / Friday lab on row 13)
$ws.Range("E13").Value = "HT LAB"
$ws.Range("F13").Value = "TOM LAB"

# Row 19: Tuesday changes from MD to FM, Friday changes from "MD - 1" to "HT"
$ws.Range("C19").Value = "FM"
$ws.Range("F19").Value = "HT"

# Row 22: Tuesday changes from FM to MD
$ws.Range("C22").Value = "MD"

# Row 25: Friday changes from HT to "MD - 1"
$ws.Range("F25").Value = "MD - 1"
